# NIT-9006246367.xlsx - "Estado de Cuenta" update
#
# The "Periodo Mora" / "Valor Mora" table (rows 16-60, columns E/F) is
# refreshed: periods are now listed in ascending chronological order
# (1607 .. 2003 instead of 2003 .. 1607) and the overdue-value (F column)
# figures attached to each period are updated to match the new database
# extract ("parte 1 de nuevos estado de cuenta").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New chronological period labels for E16:E60 (ascending 1607 -> 2003,
# replacing the previous descending 2003 -> 1607 ordering).
$periods = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

# New "Valor Mora" amounts for F16:F60, aligned 1-for-1 with $periods above.
$valores = @(
    12320,12320,12320,12320,12320,12320,
    18480,18480,18480,18480,18480,18480,18480,18480,18480,18480,18480,18480,
    24640,24640,24640,24640,
    18480,18480,18480,18480,
    31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
